$p = $ppt.ActivePresentation

# ----- Slide index 20 (sldId 303): renumber "Exercice 12/13/14" -> "11/12/13" -----
$s303 = $p.Slides.Item(20)
$s303.Shapes.Item(5).TextFrame.TextRange.Text  = "Exercice 11"
$s303.Shapes.Item(7).TextFrame.TextRange.Text  = "Exercice 12"
$s303.Shapes.Item(10).TextFrame.TextRange.Text = "Exercice 13"

# ----- Slide index 21 (sldId 304): renumber "Exercice 15/16" -> "14/15" and fix references -----
$s304 = $p.Slides.Item(21)

# "Exercice 15" -> "Exercice 14"
$s304.Shapes.Item(5).TextFrame.TextRange.Text = "Exercice 14"

# "Reprendre l'algorithme de l'exercice 14 ..." -> "... l'exercice 13 ..." (stays a single run)
$s304.Shapes.Item(6).TextFrame.TextRange.Characters(1, 122).Text = "Reprendre l’algorithme de l’exercice 13 mais cette fois-ci il faut afficher la position à laquelle le nombre a été saisi. "

# "Exercice 16" -> "Exercice 15"
$s304.Shapes.Item(7).TextFrame.TextRange.Text = "Exercice 15"

# Last paragraph: "Reprendre l'exercice 15 mais ..." becomes three runs referencing "l'exercice 14"
$tr = $s304.Shapes.Item(8).TextFrame.TextRange
$tr.Text = "Reprendre l’exercice 14 mais cette fois-ci demander à l’utilisateur la quantité de nombres qu’il veut saisir et ensuite exécuter le programme en fonction de la quantité de nombres que l’utilisateur veut saisir. Attention aux fausses saisies sur la quantité par exemple s’il saisit un nombre négatif ou 0 afficher une erreur et redemander combien de nombres il veut saisir."
$tr.Characters(11, 14).Text = "l’exercice 14 "
